# "Generate Report for Archive"
#
# The upstream commit re-runs the OpenLocalization report generator, which
# re-emits this workbook from the same underlying handoff/handback data.
# Diffing the regenerated OOXML against the previous copy shows changes
# confined to xl/sharedStrings.xml (a new, currently-unused "In Translation"
# status label is interned into the shared-string table, and the
# already-present "Ready for handoff" label is re-interned to a later slot)
# and to the <v> shared-string index references in xl/worksheets/sheet1.xml,
# sheet2.xml and sheet3.xml that track that re-indexing.
#
# Resolving every one of those index changes against the shared-string
# table (old index -> old text -> new index -> new text) shows every single
# cell keeps EXACTLY the text it had before:
#   - "Overview" (sheet1)  : A1:G5 unchanged
#   - "zh-cn"    (sheet2)  : A1:P5 unchanged
#   - "de-de"    (sheet3)  : A1:P5 unchanged
# i.e. the regeneration is a content-stable re-save: same file names, same
# paths, same statuses ("Handed back: in sync with en-US" / "Ready for
# handoff"), same handoff/handback timestamps, same xliff file names, same
# flags. Nothing a reader of the workbook (or any cell-level model) can
# observe actually changed - only the generator's internal string-interning
# order shifted, which is not something exposed through the Excel object
# model (shared-string slot order is an implementation detail the host
# manages for you from live cell contents, not something script code sets
# directly).
#
# So "applying" this edit through Excel COM automation means reproducing
# the regenerated report faithfully: touch every sheet (as the archival
# regeneration pipeline does when it re-emits the workbook) while leaving
# every cell's value exactly as it already is.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Re-generating the report revisits each sheet, but every cell already
    # holds the correct, current data - so there is nothing to overwrite.
    $ws.Calculate()
}

$excel.CalculateFullRebuild()
